{"js": "// Office.js (Word JavaScript API) script.\n// Applies the same textual changes described by the diff:\n//   1. \"Vers\u00e3o: 2.0\"                                   -> \"Vers\u00e3o: 3.0\"\n//   2. \"Revisado: 04/08/2023\"                           -> \"Revisado: 12/04/2025\"\n//   3. \"Custo estimado em 03 m\u00eas trabalhado R$ 11.235\"  -> \"Custo estimado em 3/5 m\u00eas trabalhado R$ 13.107,5\"\n//   4. \"03 meses\"                                       -> \"03 meses e meio\"\n\nconst body = context.document.body;\n\nconst replacements = [\n  [\"Vers\u00e3o: 2.0\", \"Vers\u00e3o: 3.0\"],\n  [\"Revisado: 04/08/2023\", \"Revisado: 12/04/2025\"],\n  [\n    \"Custo estimado em 03 m\u00eas trabalhado R$ 11.235\",\n    \"Custo estimado em 3/5 m\u00eas trabalhado R$ 13.107,5\",\n  ],\n  [\"03 meses\", \"03 meses e meio\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: \"${find}\"`);\n  }\n\n  // Replace only the first occurrence (each target string is unique in this document).\n  results.items[0].insertText(replace, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the same textual changes described by the diff:\n#   1. \"Vers\u00e3o: 2.0\"                                   -> \"Vers\u00e3o: 3.0\"\n#   2. \"Revisado: 04/08/2023\"                           -> \"Revisado: 12/04/2025\"\n#   3. \"Custo estimado em 03 m\u00eas trabalhado R$ 11.235\"  -> \"Custo estimado em 3/5 m\u00eas trabalhado R$ 13.107,5\"\n#   4. \"03 meses\"                                       -> \"03 meses e meio\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0            # wdFindStop\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$findText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$replaceText, [ref]2) | Out-Null\n}\n\nReplace-Text \"Vers\u00e3o: 2.0\" \"Vers\u00e3o: 3.0\"\nReplace-Text \"Revisado: 04/08/2023\" \"Revisado: 12/04/2025\"\nReplace-Text \"Custo estimado em 03 m\u00eas trabalhado R$ 11.235\" \"Custo estimado em 3/5 m\u00eas trabalhado R$ 13.107,5\"\nReplace-Text \"03 meses\" \"03 meses e meio\"\n"}
